$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 104 & 105: append two new monthly data rows -----------------

# 1) Copy the date number-format/style from the last existing row (103)
#    down into the new A104:A105 cells so they keep style index "1"
#    (the same date format used throughout column A).
$ws.Range("A103").Copy() | Out-Null
$ws.Range("A104:A105").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Columns D, E, F, G hold numeric-looking values that are stored as
#    TEXT (shared strings) in the source data, not as numbers. Force
#    those cells to text before assigning so Excel doesn't reinterpret
#    them as numbers, then clear the temporary text format again so the
#    cells end up with the default (unstyled) look, matching the rest
#    of the sheet.
$textRange = $ws.Range("D104:G105")
$textRange.NumberFormat = "@"

$ws.Cells.Item(104, 1).Value = 45839
$ws.Cells.Item(104, 2).Value = 135.039621870211
$ws.Cells.Item(104, 3).Value = 124.457525045089
$ws.Cells.Item(104, 4).Value = "116.4"
$ws.Cells.Item(104, 5).Value = "116.5"
$ws.Cells.Item(104, 6).Value = " 86.9"
$ws.Cells.Item(104, 7).Value = "173.6"

$ws.Cells.Item(105, 1).Value = 45870
$ws.Cells.Item(105, 2).Value = 135.084445598152
$ws.Cells.Item(105, 3).Value = 125.010473419394
$ws.Cells.Item(105, 4).Value = "116.8"
$ws.Cells.Item(105, 5).Value = "117.0"
$ws.Cells.Item(105, 6).Value = " 88.1"
$ws.Cells.Item(105, 7).Value = "174.2"

$textRange.ClearFormats()
